# Auto-generated: update market-price-derived columns (H-N) on several
# sheets to reflect a refreshed Universalis price pull, per the scheduled
# runner commit. Only value cells change; labels/formulas are untouched.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 3061
$ws.Range("J80").Value = 6666.3335
$ws.Range("L80").Value = 19999.0005
$ws.Range("N80").Value = -21995.0005
$ws.Range("H83").Value = 3061
$ws.Range("J83").Value = 6666.3335
$ws.Range("L83").Value = 59997.0015
$ws.Range("N83").Value = -69981.0015
$ws.Range("H88").Value = 1849.619
$ws.Range("I88").Value = 595.6667
$ws.Range("J88").Value = 2058.611
$ws.Range("K88").Value = 595.6667
$ws.Range("L88").Value = 2058.611
$ws.Range("M88").Value = -189.6667
$ws.Range("N88").Value = -2870.611
$ws.Range("H91").Value = 1849.619
$ws.Range("I91").Value = 595.6667
$ws.Range("J91").Value = 2058.611
$ws.Range("K91").Value = 595.6667
$ws.Range("L91").Value = 2058.611
$ws.Range("M91").Value = 808.3333
$ws.Range("N91").Value = -4866.611
$ws.Range("H101").Value = 3433.75
$ws.Range("I101").Value = 3911.6667
$ws.Range("J101").Value = 2000
$ws.Range("K101").Value = 11735.0001
$ws.Range("L101").Value = 6000
$ws.Range("M101").Value = -10113.0001
$ws.Range("N101").Value = -9244
$ws.Range("H107").Value = 718.1111
$ws.Range("I107").Value = 689.75
$ws.Range("K107").Value = 689.75
$ws.Range("M107").Value = 1230.25
$ws.Range("H132").Value = 4666.5
$ws.Range("I132").Value = 1229.1428
$ws.Range("K132").Value = 3687.4284
$ws.Range("M132").Value = -1157.4284
$ws.Range("H138").Value = 3976.5881
$ws.Range("I138").Value = 5665.3335
$ws.Range("J138").Value = 3614.7144
$ws.Range("K138").Value = 16996.0005
$ws.Range("L138").Value = 10844.1432
$ws.Range("M138").Value = -11856.0005
$ws.Range("N138").Value = -21124.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 602.0526
$ws.Range("I74").Value = 602.0526
$ws.Range("K74").Value = 602.0526
$ws.Range("M74").Value = 271.9474
$ws.Range("H77").Value = 602.0526
$ws.Range("I77").Value = 602.0526
$ws.Range("K77").Value = 3010.263
$ws.Range("M77").Value = 1357.737
$ws.Range("H122").Value = 2248.8333
$ws.Range("I122").Value = 2236
$ws.Range("K122").Value = 6708
$ws.Range("M122").Value = -4258
$ws.Range("H132").Value = 3070
$ws.Range("I132").Value = 3070
$ws.Range("K132").Value = 9210
$ws.Range("M132").Value = -6680

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1004.7222
$ws.Range("I107").Value = 1046.2354
$ws.Range("J107").Value = 299
$ws.Range("K107").Value = 1046.2354
$ws.Range("L107").Value = 299
$ws.Range("M107").Value = 873.7646
$ws.Range("N107").Value = -4139

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5143
$ws.Range("I31").Value = 3543.375
$ws.Range("J31").Value = 7702.4
$ws.Range("K31").Value = 3543.375
$ws.Range("L31").Value = 7702.4
$ws.Range("M31").Value = -3248.375
$ws.Range("N31").Value = -8292.4
$ws.Range("H34").Value = 5143
$ws.Range("I34").Value = 3543.375
$ws.Range("J34").Value = 7702.4
$ws.Range("K34").Value = 3543.375
$ws.Range("L34").Value = 7702.4
$ws.Range("M34").Value = -3341.375
$ws.Range("N34").Value = -8106.4
$ws.Range("H58").Value = 5078.0557
$ws.Range("I58").Value = 4239.75
$ws.Range("K58").Value = 4239.75
$ws.Range("M58").Value = -4036.75
$ws.Range("H134").Value = 1675.2222
$ws.Range("I134").Value = 1692.4546
$ws.Range("J134").Value = 1599.4
$ws.Range("K134").Value = 5077.3638
$ws.Range("L134").Value = 4798.200000000001
$ws.Range("M134").Value = -2542.3638
$ws.Range("N134").Value = -9868.200000000001
$ws.Range("H136").Value = 5078.0557
$ws.Range("I136").Value = 4239.75
$ws.Range("K136").Value = 12719.25
$ws.Range("M136").Value = -10169.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 2288.875
$ws.Range("J12").Value = 3104.4
$ws.Range("L12").Value = 9313.200000000001
$ws.Range("N12").Value = -9659.200000000001
$ws.Range("H39").Value = 252
$ws.Range("J39").Value = 252
$ws.Range("L39").Value = 756
$ws.Range("N39").Value = -1344
$ws.Range("H55").Value = 1057.1428
$ws.Range("I55").Value = 400
$ws.Range("J55").Value = 1166.6666
$ws.Range("K55").Value = 1200
$ws.Range("L55").Value = 3499.9998
$ws.Range("M55").Value = -1023
$ws.Range("N55").Value = -3853.9998
$ws.Range("H68").Value = 1649.75
$ws.Range("I68").Value = 1649.75
$ws.Range("K68").Value = 4949.25
$ws.Range("M68").Value = -4138.25
$ws.Range("H71").Value = 1649.75
$ws.Range("I71").Value = 1649.75
$ws.Range("K71").Value = 14847.75
$ws.Range("M71").Value = -10791.75
$ws.Range("H75").Value = 246.5
$ws.Range("I75").Value = 393
$ws.Range("J75").Value = 100
$ws.Range("K75").Value = 1179
$ws.Range("L75").Value = 300
$ws.Range("M75").Value = -181
$ws.Range("N75").Value = -2296
$ws.Range("H78").Value = 246.5
$ws.Range("I78").Value = 393
$ws.Range("J78").Value = 100
$ws.Range("K78").Value = 3537
$ws.Range("L78").Value = 900
$ws.Range("M78").Value = 1455
$ws.Range("N78").Value = -10884
$ws.Range("H107").Value = 965
$ws.Range("I107").Value = 394.4
$ws.Range("J107").Value = 1168.7858
$ws.Range("K107").Value = 1183.2
$ws.Range("L107").Value = 3506.3574
$ws.Range("M107").Value = 736.8000000000002
$ws.Range("N107").Value = -7346.357400000001
$ws.Range("H113").Value = 261.16666
$ws.Range("I113").Value = 248.5
$ws.Range("J113").Value = 267.5
$ws.Range("K113").Value = 745.5
$ws.Range("L113").Value = 802.5
$ws.Range("M113").Value = 1424.5
$ws.Range("N113").Value = -5142.5
$ws.Range("H132").Value = 3500
$ws.Range("J132").Value = 3500
$ws.Range("L132").Value = 31500
$ws.Range("N132").Value = -36560
$ws.Range("H137").Value = 4500
$ws.Range("I137").Value = 4500
$ws.Range("K137").Value = 13500
$ws.Range("M137").Value = -8400
$ws.Range("H138").Value = 1500.5
$ws.Range("I138").Value = 1000
$ws.Range("K138").Value = 3000
$ws.Range("M138").Value = 2140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H122").Value = 5056.6665
$ws.Range("I122").Value = 4743.4443
$ws.Range("J122").Value = 5996.3335
$ws.Range("K122").Value = 14230.3329
$ws.Range("L122").Value = 17989.0005
$ws.Range("M122").Value = -11780.3329
$ws.Range("N122").Value = -22889.0005
$ws.Range("H132").Value = 2821.8333
$ws.Range("I132").Value = 2821.8333
$ws.Range("K132").Value = 8465.499899999999
$ws.Range("M132").Value = -5935.499899999999
$ws.Range("H136").Value = 4001.5
$ws.Range("I136").Value = 3751
$ws.Range("J136").Value = 4502.5
$ws.Range("K136").Value = 11253
$ws.Range("L136").Value = 13507.5
$ws.Range("M136").Value = -8703
$ws.Range("N136").Value = -18607.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 812.5
$ws.Range("I81").Value = 750
$ws.Range("K81").Value = 1500
$ws.Range("M81").Value = -439
$ws.Range("H84").Value = 812.5
$ws.Range("I84").Value = 750
$ws.Range("K84").Value = 7500
$ws.Range("M84").Value = -2196
$ws.Range("H105").Value = 19999.5
$ws.Range("J105").Value = 19999.5
$ws.Range("L105").Value = 19999.5
$ws.Range("N105").Value = -26987.5
$ws.Range("H132").Value = 39751.848
$ws.Range("I132").Value = 44719.914
$ws.Range("K132").Value = 134159.742
$ws.Range("M132").Value = -131629.742
$ws.Range("H136").Value = 3186.25
$ws.Range("I136").Value = 1648.3334
$ws.Range("J136").Value = 7800
$ws.Range("K136").Value = 4945.0002
$ws.Range("L136").Value = 23400
$ws.Range("M136").Value = -2395.0002
$ws.Range("N136").Value = -28500
